$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3, shifting the table down by one row
$ws.Rows("3:3").Insert()

# Make row 2 (the "Phrase 1" header) bold to match the new style
$ws.Range("A2").Font.Bold = $true

# Update the selection to match the target state
$ws.Range("G6").Select()
